$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (shifts RuleTable/Condition/Action row, U/date/customer name row,
# and 1/foo/cust.amount.value row each down by one)
$ws.Rows.Item(3).Insert()

# Fill the newly inserted row 3 with "variable" / "foo" (quoted)
$ws.Cells.Item(3, 2).Value = "variable"
$ws.Cells.Item(3, 3).Value = '"foo"'

# Row 5 (originally the U / date / customer name row, now shifted down to row 5):
# "date" becomes "gate"
$ws.Cells.Item(5, 3).Value = "gate"

# Row 6 (originally the 1 / foo / cust.amount.value row, now shifted down to row 6):
# "foo" becomes "foo" (quoted) and cust.amount.value becomes numeric 100
$ws.Cells.Item(6, 3).Value = '"foo"'
$ws.Cells.Item(6, 4).Value = 100

# Update the active selection to C7
$ws.Range("C7").Select()
